$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PID to PD switch: desired_normal_force (column C) changes from 1 to 0.5
# and sliding_velocity (column D) changes from 20 to 10 for row 3.
# Rows 4-6 are newly populated with the same pattern as rows 2-3.

$ws.Range("C2").Value = 0.5

$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 10

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 10

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 10

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 10

# Update the active selection to G10
$ws.Range("G10").Select()
